$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CompStat_1")

$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F14").PasteSpecial(-4122)

$ws.Range("C15").Value = 1
$ws.Range("G15").Copy()
$ws.Range("C15").PasteSpecial(-4122)

$ws.Range("D15").Value = 1
$ws.Range("G15").Copy()
$ws.Range("D15").PasteSpecial(-4122)

$ws.Range("E15").Value = 0
$ws.Range("H15").Copy()
$ws.Range("E15").PasteSpecial(-4122)

$ws.Range("F15").Value = 1
$ws.Range("G15").Copy()
$ws.Range("F15").PasteSpecial(-4122)

$ws.Range("G15").Value = 4

$ws.Range("H15").Value = -75

$ws.Range("I15").Value = 1
$ws.Range("G15").Copy()
$ws.Range("I15").PasteSpecial(-4122)

$ws.Range("J15").Value = 1
$ws.Range("G15").Copy()
$ws.Range("J15").PasteSpecial(-4122)

$ws.Range("K15").Value = 0
$ws.Range("H15").Copy()
$ws.Range("K15").PasteSpecial(-4122)

$ws.Range("D16").Value = 1

$ws.Range("E16").Value = 300

$ws.Range("F16").Value = 16

$ws.Range("G16").Value = 11

$ws.Range("H16").Value = 45.454545454545

$ws.Range("I16").Value = 6

$ws.Range("J16").Value = 1
$ws.Range("G15").Copy()
$ws.Range("J16").PasteSpecial(-4122)

$ws.Range("K16").Value = 500
$ws.Range("H15").Copy()
$ws.Range("K16").PasteSpecial(-4122)

$ws.Range("L16").Value = 50
$ws.Range("H15").Copy()
$ws.Range("L16").PasteSpecial(-4122)

$ws.Range("M16").Value = 500
$ws.Range("H15").Copy()
$ws.Range("M16").PasteSpecial(-4122)

$ws.Range("N16").Value = -76.923076923076

$ws.Range("C17").Value = 3

$ws.Range("E17").Value = 0

$ws.Range("F17").Value = 12

$ws.Range("G17").Value = 11

$ws.Range("H17").Value = 9.090909090909

$ws.Range("I17").Value = 5

$ws.Range("J17").Value = 3
$ws.Range("G15").Copy()
$ws.Range("J17").PasteSpecial(-4122)

$ws.Range("K17").Value = 66.666666666666
$ws.Range("H15").Copy()
$ws.Range("K17").PasteSpecial(-4122)

$ws.Range("L17").Value = 400
$ws.Range("H15").Copy()
$ws.Range("L17").PasteSpecial(-4122)

$ws.Range("M17").Value = 400
$ws.Range("H15").Copy()
$ws.Range("M17").PasteSpecial(-4122)

$ws.Range("N17").Value = -44.444444444444

$ws.Range("C18").Value = 9

$ws.Range("D18").Value = 6

$ws.Range("E18").Value = 50

$ws.Range("F18").Value = 28

$ws.Range("G18").Value = 20

$ws.Range("H18").Value = 40

$ws.Range("I18").Value = 10

$ws.Range("J18").Value = 9

$ws.Range("K18").Value = 11.111111111111

$ws.Range("L18").Value = 11.111111111111
$ws.Range("H15").Copy()
$ws.Range("L18").PasteSpecial(-4122)

$ws.Range("M18").Value = 233.333333333333
$ws.Range("H15").Copy()
$ws.Range("M18").PasteSpecial(-4122)

$ws.Range("N18").Value = -28.571428571428

$ws.Range("C19").Value = 30

$ws.Range("D19").Value = 21

$ws.Range("E19").Value = 42.857142857142

$ws.Range("F19").Value = 97

$ws.Range("G19").Value = 82

$ws.Range("H19").Value = 18.292682926829

$ws.Range("I19").Value = 35

$ws.Range("J19").Value = 25

$ws.Range("K19").Value = 40

$ws.Range("L19").Value = 337.5
$ws.Range("H15").Copy()
$ws.Range("L19").PasteSpecial(-4122)

$ws.Range("M19").Value = 66.666666666666

$ws.Range("N19").Value = -30

$ws.Range("C20").Value = 1
$ws.Range("G15").Copy()
$ws.Range("C20").PasteSpecial(-4122)

$ws.Range("F20").Value = 2

$ws.Range("G20").Value = 1

$ws.Range("H20").Value = 100

$ws.Range("I20").Value = 1
$ws.Range("G15").Copy()
$ws.Range("I20").PasteSpecial(-4122)

$ws.Range("L20").Value = 0
$ws.Range("H15").Copy()
$ws.Range("L20").PasteSpecial(-4122)

$ws.Range("N20").Value = -92.857142857142

$ws.Range("C21").Value = 48

$ws.Range("D21").Value = 32

$ws.Range("E21").Value = 50

$ws.Range("F21").Value = 156

$ws.Range("G21").Value = 129

$ws.Range("H21").Value = 20.930232558139

$ws.Range("I21").Value = 58

$ws.Range("J21").Value = 39

$ws.Range("K21").Value = 48.717948717948

$ws.Range("L21").Value = 152.173913043478
$ws.Range("E21").Copy()
$ws.Range("L21").PasteSpecial(-4122)

$ws.Range("M21").Value = 123.076923076923

$ws.Range("N21").Value = -48.672566371681

$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C22").PasteSpecial(-4122)

$ws.Range("E22").Value = -100

$ws.Range("F22").Value = 5

$ws.Range("G22").Value = 3

$ws.Range("H22").Value = 66.666666666666

$ws.Range("J22").Value = 1
$ws.Range("G15").Copy()
$ws.Range("J22").PasteSpecial(-4122)

$ws.Range("K22").Value = -100
$ws.Range("H15").Copy()
$ws.Range("K22").PasteSpecial(-4122)

$ws.Range("L22").Value = -100
$ws.Range("H15").Copy()
$ws.Range("L22").PasteSpecial(-4122)

$ws.Range("M22").Value = -100
$ws.Range("H15").Copy()
$ws.Range("M22").PasteSpecial(-4122)

$ws.Range("C24").Value = 35

$ws.Range("D24").Value = 29

$ws.Range("E24").Value = 20.689655172413

$ws.Range("G24").Value = 112

$ws.Range("H24").Value = 11.607142857142

$ws.Range("I24").Value = 37

$ws.Range("J24").Value = 30

$ws.Range("K24").Value = 23.333333333333

$ws.Range("L24").Value = 2.777777777777

$ws.Range("M24").Value = 48
$ws.Range("H15").Copy()
$ws.Range("M24").PasteSpecial(-4122)

$ws.Range("C25").Value = 11

$ws.Range("E25").Value = 266.666666666667

$ws.Range("F25").Value = 28

$ws.Range("G25").Value = 27

$ws.Range("H25").Value = 3.703703703703

$ws.Range("I25").Value = 14

$ws.Range("J25").Value = 3
$ws.Range("G15").Copy()
$ws.Range("J25").PasteSpecial(-4122)

$ws.Range("K25").Value = 366.666666666667
$ws.Range("H15").Copy()
$ws.Range("K25").PasteSpecial(-4122)

$ws.Range("L25").Value = 133.333333333333
$ws.Range("H15").Copy()
$ws.Range("L25").PasteSpecial(-4122)

$ws.Range("M25").Value = 366.666666666667

$ws.Range("C26").Value = 1
$ws.Range("G15").Copy()
$ws.Range("C26").PasteSpecial(-4122)

$ws.Range("D26").Value = 1
$ws.Range("G15").Copy()
$ws.Range("D26").PasteSpecial(-4122)

$ws.Range("E26").Value = 0
$ws.Range("H15").Copy()
$ws.Range("E26").PasteSpecial(-4122)

$ws.Range("F26").Value = 1
$ws.Range("G15").Copy()
$ws.Range("F26").PasteSpecial(-4122)

$ws.Range("G26").Value = 4

$ws.Range("H26").Value = -75

$ws.Range("I26").Value = 1
$ws.Range("G15").Copy()
$ws.Range("I26").PasteSpecial(-4122)

$ws.Range("J26").Value = 1
$ws.Range("G15").Copy()
$ws.Range("J26").PasteSpecial(-4122)

$ws.Range("K26").Value = 0
$ws.Range("H15").Copy()
$ws.Range("K26").PasteSpecial(-4122)

$ws.Range("F27").Value = 4

$ws.Range("H27").Value = 100

$ws.Range("I27").Value = 2

$ws.Range("J27").Value = 1
$ws.Range("G15").Copy()
$ws.Range("J27").PasteSpecial(-4122)

$ws.Range("K27").Value = 100
$ws.Range("H15").Copy()
$ws.Range("K27").PasteSpecial(-4122)

$ws.Range("L27").Value = 100
$ws.Range("H15").Copy()
$ws.Range("L27").PasteSpecial(-4122)

$ws.Range("J40").Value = 424

$ws.Range("K40").Value = 44.217687074829

$ws.Range("L40").Value = 11.578947368421

$ws.Range("M40").Value = -47.654320987654

$ws.Range("N40").Value = -64.128595600676

$ws.Range("J41").Value = 1372

$ws.Range("K41").Value = 12.274959083469

$ws.Range("L41").Value = -15.413070283600

$ws.Range("M41").Value = -46.697746697746

$ws.Range("N41").Value = -64.224250325945

$ws.Range("J43").Value = 2297

$ws.Range("K43").Value = 8.349056603773

$ws.Range("L43").Value = -22.029871011541

$ws.Range("M43").Value = -56.968902210565

$ws.Range("N43").Value = -70.694054605766

# Shared-string partial text edits (preserve rich-text run formatting)
$a8 = $ws.Range("A8")
$a8.Characters(21,2).Text = "1"

$c9 = $ws.Range("C9")
$c9.Characters(27,10).Text = "1/2/2023"
$c9.Characters(46,8).Text = "1/8/2023"

$excel.CutCopyMode = $false